# componentes.xlsx v0.0.1 -> v0.0.2 component list update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update component names / quantities (rows 2-13) ---
$ws.Range("A2").Value = "ESP32"

$ws.Range("A4").Value = "Buzzer 5V 12mm"

$ws.Range("A6").Value = "Interruptor"

$ws.Range("A7").Value = "Jumper MF"
$ws.Range("B7").Value = "~30"

$ws.Range("A8").Value = "Jumper MM"
$ws.Range("B8").Value = "~10"

$ws.Range("A10").Value = "Led Vermelho"
$ws.Range("B10").Value = 4

$ws.Range("A13").Value = "Suporte Para Bateria 4 Slots"

$ws.Range("A9").Value = "Led Amarelo"
$ws.Range("B9").Value = 2

# Row 7 (now "Jumper MF") previously carried the one-off highlight style
# used by the old "Interruptor" row; bring it back in line with the
# regular row styling (copy formatting from a normal row).
$ws.Range("A6:B6").Copy()
$ws.Range("A7:B7").PasteSpecial(-4122)

# --- Update the saved cursor/selection position ---
$ws.Range("A17").Select()
